# add MoCo - w/o distributed yet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New experiment row (row 6): fill in best_epoch_id + accuracy columns ---
$ws.Range("F6").Value = 16

$ws.Range("J6").Value = 0.66279999999999994
$ws.Range("K6").Value = 0.58640000000000003
$ws.Range("L6").Value = 0.59770000000000001

$ws.Range("J6:L6").NumberFormat = "0.00%"

# --- Column width adjustments ---
# col J (10): widened to fit new numbers, now explicitly custom-sized
$ws.Columns.Item(10).ColumnWidth = 14.498697916666666
# col L (12): split off from the 12:35 default-width run, given its own width
$ws.Columns.Item(12).ColumnWidth = 20.166666666666668

# --- View state: scrolled right to show the accuracy columns, zoomed in, new selection ---
$ws.Range("L5").Select()
$excel.ActiveWindow.Zoom = 157
